# Generate Report for Handback
# Updates the handoff/handback generation timestamps for the
# "8e5b40ae-ac50-430c-b4df-cb1eead22b23" entry (row 2) across all report
# sheets, reflecting a freshly generated handback report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# 8e5b40ae-ac50-430c-b4df-cb1eead22b23 row.
$wsOverview.Range("G2").Value = "2016-08-15 12:44:10"

# zh-cn sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) columns for the same row.
$wsZhCn.Range("H2").Value = "2016-08-15 12:44:02"
$wsZhCn.Range("K2").Value = "2016-08-15 12:44:29"

# de-de sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) columns for the same row.
$wsDeDe.Range("H2").Value = "2016-08-15 12:44:10"
$wsDeDe.Range("K2").Value = "2016-08-15 12:44:35"
